# Update SwaadSutra_Consolidated_2026-01-13.xlsx
# A brand-new order (Order ID 13, "Ketki") was placed, so a new row is
# inserted at the top of the order list (row 2) on the "All Orders" sheet,
# pushing every existing order down by one row. The "Daily Summary" sheet
# totals are updated accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)       # "All Orders"
$summary = $wb.Worksheets.Item(2)  # "Daily Summary"

# ---------------------------------------------------------------------
# 1. Insert a new row right under the header row, shifting all existing
#    orders down by one row.
# ---------------------------------------------------------------------
$ws.Range("A2").EntireRow.Insert()

# ---------------------------------------------------------------------
# 2. Populate the newly inserted row with the new order's data.
#    Columns D, E and J look like plain numbers / dates, so force them to
#    be stored as text (matching the rest of the sheet) before assigning.
# ---------------------------------------------------------------------
$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("J2").NumberFormat = "@"

$ws.Range("A2").Value2 = 13
$ws.Range("B2").Value2 = "2026-01-13 22:51"
$ws.Range("C2").Value2 = "Ketki"
$ws.Range("D2").Value2 = "1608"
$ws.Range("E2").Value2 = "3159135521"
$ws.Range("F2").Value2 = "Wheat Chapati x1"
$ws.Range("G2").Value2 = 15
$ws.Range("H2").Value2 = "NEW"
$ws.Range("I2").Value2 = "PENDING"
$ws.Range("J2").Value2 = "2026-01-14"
$ws.Range("K2").Value2 = "16:51"
# L2, M2 and N2 (Notes / Cancel Reason / Feedback) are left blank for the
# new order, matching the blank cells used throughout the rest of the sheet.

# ---------------------------------------------------------------------
# 3. Update the "Daily Summary" sheet totals for 2026-01-13.
# ---------------------------------------------------------------------
$summary.Range("B2").Value2 = 13   # Total Orders
$summary.Range("E2").Value2 = 340  # Revenue
$summary.Range("G2").Value2 = 340  # Pending
